# Product card generation update
# Regenerate the "Характеристики" (specs) text for the С20 profiled-sheet
# row (row 5, column H on sheet "Кровля") using the updated card template:
# each "Label" / "Value" pair becomes "Label:" / "Value<br>" joined by line
# breaks, matching the new HTML product-card generator output.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Кровля")

$nl = [char]10

$pairs = @(
    @("Тип продукта", "Профнастил"),
    @("Основной материал", "Сталь"),
    @("Цветовая палитра", "Серый / серебристый"),
    @("Цветовая палитра по RAL", "Нет"),
    @("Площадь покрытия продуктом (м²)", "2.32"),
    @("Полезная площадь (м²)", "2.09"),
    @("Длина (см)", "200"),
    @("Ширина (см)", "116"),
    @("Толщина (мм)", "0.35"),
    @("Высота волны (мм)", "20"),
    @("Вес на м² (кг)", "3,1"),
    @("Вес, кг", "6.5"),
    @("Страна производства", "Россия")
)

$parts = @()
for ($i = 0; $i -lt $pairs.Count; $i++) {
    $label = $pairs[$i][0]
    $value = $pairs[$i][1]
    $parts += ($label + ":")
    if ($i -eq ($pairs.Count - 1)) {
        $parts += $value
    } else {
        $parts += ($value + "<br>")
    }
}

$newCharacteristics = $parts -join $nl

$ws.Range("H5").Value = $newCharacteristics

# Assigning multi-line content makes Excel auto-fit the row; restore the
# original fixed row height that the sheet used before the edit.
$ws.Rows.Item(5).RowHeight = 17.25

# Update the saved cursor/selection position on the active sheet.
$ws.Activate()
$ws.Range("E14").Select()
